# "Generate Report for Handback"
#
# The localization status report is refreshed: the "Ready for handoff"
# status becomes "Handed back: in sync with en-US" everywhere it appears
# (Overview sheet + each per-language sheet), the Latest Handback
# timestamps for zh-cn/de-de are bumped to the new handback run, and the
# previously-recorded "handback file is not latest" error detail is
# cleared now that the handback is in sync.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status cells -----------------------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# Match the column auto-fit that happens after the longer status text
# is written (closest value this host's ColumnWidth quantization can
# reach to the real-Excel autofit width).
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-09-05 09:08:40"
$wsZhCn.Range("P2").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

# --- de-de sheet -------------------------------------------------------
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-09-05 09:08:48"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
